$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text edits (rich-text runs) ---
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "18"

$c9 = $ws.Range("C9")
$c9.Characters(27, 9).Text = "4/29/2024"
$c9.Characters(47, 9).Text = "5/5/2024"

# --- Data table value edits (rows 15-30) ---
# Row 15
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 7
$ws.Range("H15").Value = -71.428571428571
$ws.Range("I15").Value = 16
$ws.Range("J15").Value = 16
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 14.285714285714
$ws.Range("M15").Value = 220
$ws.Range("N15").Value = -23.809523809523

# Row 16
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = 25
$ws.Range("F16").Value = 38
$ws.Range("G16").Value = 31
$ws.Range("H16").Value = 22.58064516129
$ws.Range("I16").Value = 145
$ws.Range("J16").Value = 130
$ws.Range("K16").Value = 11.538461538461
$ws.Range("L16").Value = 33.027522935779
$ws.Range("M16").Value = 74.698795180722
$ws.Range("N16").Value = -63.659147869674

# Row 17
$ws.Range("C17").Value = 17
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 70
$ws.Range("F17").Value = 63
$ws.Range("G17").Value = 45
$ws.Range("H17").Value = 40
$ws.Range("I17").Value = 241
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 20.5
$ws.Range("L17").Value = 16.990291262135
$ws.Range("M17").Value = 180.232558139535
$ws.Range("N17").Value = -19.397993311036

# Row 18
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 23
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 91.666666666666
$ws.Range("I18").Value = 98
$ws.Range("J18").Value = 71
$ws.Range("K18").Value = 38.028169014084
$ws.Range("L18").Value = -11.711711711711
$ws.Range("M18").Value = 139.024390243902
$ws.Range("N18").Value = -72.549019607843

# Row 19
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 25
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 35
$ws.Range("H19").Value = 8.571428571428
$ws.Range("I19").Value = 191
$ws.Range("J19").Value = 153
$ws.Range("K19").Value = 24.836601307189
$ws.Range("L19").Value = 40.441176470588
$ws.Range("M19").Value = 161.643835616438
$ws.Range("N19").Value = 85.436893203883

# Row 20
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -14.285714285714
$ws.Range("F20").Value = 21
$ws.Range("G20").Value = 40
$ws.Range("H20").Value = -47.5
$ws.Range("I20").Value = 83
$ws.Range("J20").Value = 160
$ws.Range("K20").Value = -48.125
$ws.Range("L20").Value = -14.432989690721
$ws.Range("M20").Value = 102.439024390244
$ws.Range("N20").Value = -58.5

# Row 21
$ws.Range("C21").Value = 48
$ws.Range("D21").Value = 36
$ws.Range("E21").Value = 33.333333333333
$ws.Range("F21").Value = 185
$ws.Range("G21").Value = 172
$ws.Range("H21").Value = 7.558139534883
$ws.Range("I21").Value = 775
$ws.Range("J21").Value = 734
$ws.Range("K21").Value = 5.58583106267
$ws.Range("L21").Value = 14.644970414201
$ws.Range("M21").Value = 132.035928143713
$ws.Range("N21").Value = -44.124008651766

# Row 22
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 4
$ws.Range("J22").Value = 5
$ws.Range("K22").Value = -20
$ws.Range("L22").Value = 33.333333333333
$ws.Range("M22").Value = -50

# Row 23
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 27
$ws.Range("G23").Value = 38
$ws.Range("H23").Value = -28.947368421052
$ws.Range("I23").Value = 131
$ws.Range("J23").Value = 144
$ws.Range("K23").Value = -9.027777777777
$ws.Range("L23").Value = 15.929203539823
$ws.Range("M23").Value = 111.290322580645

# Row 24
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = -11.111111111111
$ws.Range("F24").Value = 90
$ws.Range("G24").Value = 72
$ws.Range("H24").Value = 25
$ws.Range("I24").Value = 386
$ws.Range("J24").Value = 381
$ws.Range("K24").Value = 1.312335958005
$ws.Range("L24").Value = 3.208556149732
$ws.Range("M24").Value = 54.4

# Row 25
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 16
$ws.Range("G25").Value = 14
$ws.Range("H25").Value = 14.285714285714
$ws.Range("I25").Value = 58
$ws.Range("J25").Value = 79
$ws.Range("K25").Value = -26.582278481012
$ws.Range("L25").Value = -46.788990825688

# Row 26
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 17
$ws.Range("E26").Value = -35.294117647058
$ws.Range("F26").Value = 65
$ws.Range("G26").Value = 83
$ws.Range("H26").Value = -21.686746987951
$ws.Range("I26").Value = 289
$ws.Range("J26").Value = 388
$ws.Range("K26").Value = -25.515463917525
$ws.Range("L26").Value = -14.749262536873
$ws.Range("M26").Value = 2.120141342756

# Row 27
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("I27").Value = 22
$ws.Range("J27").Value = 23
$ws.Range("K27").Value = -4.347826086956
$ws.Range("L27").Value = -4.347826086956

# Row 28
$ws.Range("F28").Value = 9
$ws.Range("H28").Value = 12.5
$ws.Range("I28").Value = 27
$ws.Range("J28").Value = 36
$ws.Range("K28").Value = -25
$ws.Range("L28").Value = 12.5

# Row 29
$ws.Range("F29").Value = 4
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 300
$ws.Range("I29").Value = 7
$ws.Range("K29").Value = -30
$ws.Range("L29").Value = -41.666666666666
$ws.Range("M29").Value = -53.333333333333
$ws.Range("N29").Value = -74.074074074074

# Row 30
$ws.Range("F30").Value = 4
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 300
$ws.Range("I30").Value = 7
$ws.Range("K30").Value = -30
$ws.Range("L30").Value = -36.363636363636
$ws.Range("M30").Value = -46.153846153846
$ws.Range("N30").Value = -74.074074074074

# --- Row 22: convert D/E/G/H from text placeholders to numeric cells ---
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("G22").NumberFormat = "#,##0"
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H22").NumberFormat = "#,##0.0;""-""#,##0.0"

# --- Column E width bestfit-style adjustment ---
$ws.Columns.Item(5).ColumnWidth = 6.714285714285714
